# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-29 22:21:41
#
# Normalises the "Recorded By" (column G) cell text on the "Session Analysis
# Results" sheet: entries are re-ordered so "System" / "system" always sorts
# ahead of the human recorder, e.g.
#   "dnasr281@gmail.com, System"              -> "System, dnasr281@gmail.com"
#   "admin@admin.com, System"                 -> "System, admin@admin.com"
#   "backup@backdoor.com, system, System"     -> "backup@backdoor.com, System, system"
# The "backup@backdoor.com, System" (2-token) rows are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val.Split(",")
    for ($i = 0; $i -lt $parts.Count; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    $newVal = $val

    # NOTE: PowerShell's "-eq"/"-ne" are case-INSENSITIVE, which would blur
    # "System" vs "system" here, so case-sensitive comparisons use the
    # plain .NET .Equals() method instead.
    if ($parts.Count -eq 2 -and $parts[1].Equals("System") -and -not $parts[0].Equals("backup@backdoor.com")) {
        # "<recorder>, System" -> "System, <recorder>"
        $newVal = "$($parts[1]), $($parts[0])"
    }
    elseif ($parts.Count -eq 3 -and $parts[0].Equals("backup@backdoor.com") -and `
            (($parts[1].Equals("system") -and $parts[2].Equals("System")) -or `
             ($parts[1].Equals("System") -and $parts[2].Equals("system")))) {
        # "backup@backdoor.com, system, System" -> "backup@backdoor.com, System, system"
        $newVal = "$($parts[0]), $($parts[2]), $($parts[1])"
    }

    if (-not $newVal.Equals($val)) {
        $cell.Value = $newVal
    }
}
